# Natmi following Dr Hou advice
# Update ligand/receptor expression statistics on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.574538
$ws.Range("H2").Value = 1.723614
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 0.752908
$ws.Range("N2").Value = 1.505816
$ws.Range("O2").Value = 0.04239205579776523
$ws.Range("P2").Value = 0.03064402855818915
$ws.Range("Q2").Value = 0.432574256504
$ws.Range("R2").Value = 2.595445539024
$ws.Range("S2").Value = 0.04239205579776523
$ws.Range("T2").Value = 0.03064402855818915

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.574538
$ws.Range("H3").Value = 1.723614
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 12.18022033333333
$ws.Range("N3").Value = 36.540661
$ws.Range("O3").Value = 0.6858003633906682
$ws.Range("P3").Value = 0.743618781590253
$ws.Range("Q3").Value = 6.997999429872666
$ws.Range("R3").Value = 62.981994868854
$ws.Range("S3").Value = 0.6858003633906682
$ws.Range("T3").Value = 0.743618781590253

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.574538
$ws.Range("H4").Value = 1.723614
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.4715666666666666
$ws.Range("N4").Value = 1.4147
$ws.Range("O4").Value = 0.02655129238326527
$ws.Range("P4").Value = 0.02878977723790303
$ws.Range("Q4").Value = 0.2709329695333333
$ws.Range("R4").Value = 2.4383967258
$ws.Range("S4").Value = 0.02655129238326527
$ws.Range("T4").Value = 0.02878977723790303

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.574538
$ws.Range("H5").Value = 1.723614
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.3304726666666666
$ws.Range("N5").Value = 0.9914179999999999
$ws.Range("O5").Value = 0.01860707513397334
$ws.Range("P5").Value = 0.02017579937064207
$ws.Range("Q5").Value = 0.1898691049613333
$ws.Range("R5").Value = 1.708821944652
$ws.Range("S5").Value = 0.01860707513397334
$ws.Range("T5").Value = 0.02017579937064207

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.574538
$ws.Range("H6").Value = 1.723614
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6355266666666667
$ws.Range("N6").Value = 1.90658
$ws.Range("O6").Value = 0.03578296672940263
$ws.Range("P6").Value = 0.03879975506202103
$ws.Range("Q6").Value = 0.3651342200133333
$ws.Range("R6").Value = 3.28620798012
$ws.Range("S6").Value = 0.03578296672940263
$ws.Range("T6").Value = 0.03879975506202103

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.574538
$ws.Range("H7").Value = 1.723614
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 3.3898975
$ws.Range("N7").Value = 6.779795
$ws.Range("O7").Value = 0.1908662465649254
$ws.Range("P7").Value = 0.1379718581809916
$ws.Range("Q7").Value = 1.947624929855
$ws.Range("R7").Value = 11.68574957913
$ws.Range("S7").Value = 0.1908662465649254
$ws.Range("T7").Value = 0.1379718581809916
